$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D/E hold numeric- or percent-looking text that must remain
# plain text (matching the source inlineStr cells). Temporarily switching the
# NumberFormat to Text ("@") before the assignment stops Excel from silently
# re-typing the string as a number/percentage; flipping it back to "General"
# afterwards restores the original display format.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "256.59"
Set-TextValue $ws.Range("E2") "0.38%"
Set-TextValue $ws.Range("D3") "27.11"
Set-TextValue $ws.Range("E3") "-3.72%"
Set-TextValue $ws.Range("D4") "4.662"
Set-TextValue $ws.Range("E4") "-10.65%"
Set-TextValue $ws.Range("D5") "0.05901"
Set-TextValue $ws.Range("E5") "0.73%"
Set-TextValue $ws.Range("D6") "6.650"
Set-TextValue $ws.Range("E6") "-0.81%"
Set-TextValue $ws.Range("D7") "0.8678"
Set-TextValue $ws.Range("E7") "-0.30%"
Set-TextValue $ws.Range("D8") "0.9526"
Set-TextValue $ws.Range("E8") "-0.68%"
Set-TextValue $ws.Range("E9") "-0.60%"
Set-TextValue $ws.Range("D10") "0.03721"
Set-TextValue $ws.Range("E10") "7.92%"
Set-TextValue $ws.Range("D11") "0.07080"
Set-TextValue $ws.Range("E11") "-1.12%"
Set-TextValue $ws.Range("D12") "0.03207"
Set-TextValue $ws.Range("E12") "0.03%"
Set-TextValue $ws.Range("D13") "0.09261"
Set-TextValue $ws.Range("E13") "0.55%"
Set-TextValue $ws.Range("D14") "0.001547"
Set-TextValue $ws.Range("E14") "-0.49%"
Set-TextValue $ws.Range("D15") "0.0006000"
Set-TextValue $ws.Range("E15") "-1.34%"
Set-TextValue $ws.Range("D16") "0.006044"
Set-TextValue $ws.Range("E16") "1.99%"
Set-TextValue $ws.Range("E17") "0.48%"
Set-TextValue $ws.Range("E18") "-1.12%"
Set-TextValue $ws.Range("E19") "-0.12%"
Set-TextValue $ws.Range("D20") "0.3076"
Set-TextValue $ws.Range("E20") "-3.22%"
Set-TextValue $ws.Range("E21") "-1.96%"
Set-TextValue $ws.Range("E22") "8.87%"
Set-TextValue $ws.Range("D23") "0.04237"
Set-TextValue $ws.Range("E23") "1.11%"
Set-TextValue $ws.Range("D25") "0.001219"
Set-TextValue $ws.Range("E25") "-0.18%"
Set-TextValue $ws.Range("D26") "0.004279"
Set-TextValue $ws.Range("E26") "-6.28%"
Set-TextValue $ws.Range("E28") "2.38%"
Set-TextValue $ws.Range("D40") "0.03817"
Set-TextValue $ws.Range("E40") "-0.03%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1099"
Set-TextValue $ws.Range("E41") "-0.18%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D42") "0.006216"
Set-TextValue $ws.Range("E42") "10.02%"
Set-TextValue $ws.Range("D43") "0.002289"
Set-TextValue $ws.Range("E43") "-2.35%"
Set-TextValue $ws.Range("D44") "0.01157"
Set-TextValue $ws.Range("E44") "17.84%"
Set-TextValue $ws.Range("D45") "0.00005500"
Set-TextValue $ws.Range("E45") "1.61%"
Set-TextValue $ws.Range("E46") "-0.01%"
Set-TextValue $ws.Range("D47") "0.06020"
Set-TextValue $ws.Range("E47") "-33.13%"
Set-TextValue $ws.Range("E48") "7.03%"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "-0.01%"
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "-0.01%"
